# Reorder the Csets/CommName pairs in rows 4-45 of the "Commodities" sheet
# back to an earlier ordering (commit: "retracted back to earlier method to
# avoid complications"). Row 21 and row 46 are unchanged; all other rows in
# the 4-45 range get new (B, C) values representing a permutation of the
# same set of pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commodities")

$data = @(
    @(4,  "MAT", "iip_steel_oxygen"),
    @(5,  "MAT", "iip_steel_sinter"),
    @(6,  "MAT", "iip_steel_scrap"),
    @(7,  "NRG", "pri_uran"),
    @(8,  "ENV", "emi_CO2_f_x2x_neg_reusable"),
    @(9,  "NRG", "sec_heat_high"),
    @(10, "NRG", "sec_heavy_fuel_oil"),
    @(11, "ENV", "emi_CH4_f_ind"),
    @(12, "NRG", "CO2_f_pow"),
    @(13, "NRG", "pri_biomass"),
    @(14, "ENV", "[emi_CO2_f_x2x_neg_reusable"),
    @(15, "ENV", "emi_CO2_f_ind"),
    @(16, "MAT", "iip_steel_sponge_iron"),
    @(17, "NRG", "iip_coke"),
    @(18, "MAT", "iip_steel_crudesteel"),
    @(19, "NRG", "sec_heat_low"),
    @(20, "NRG", "pri_waste"),
    @(22, "NRG", "pri_crude_oil"),
    @(23, "MAT", "iip_steel_raw_iron"),
    @(24, "NRG", "sec_biogas"),
    @(25, "NRG", "pri_hydro_energy"),
    @(26, "NRG", "sec_elec_ind"),
    @(27, "NRG", "sec_natural_gas_syn"),
    @(28, "NRG", "pri_geoth_heat"),
    @(29, "NRG", "sec_H2"),
    @(30, "NRG", "iip_heat_proc"),
    @(31, "MAT", "iip_steel_iron_pellets"),
    @(32, "NRG", "sec_elec"),
    @(33, "NRG", "pri_natural_gas"),
    @(34, "ENV", "emi_CO2_f_x2x_neg_stored]"),
    @(35, "DEM", "exo_steel"),
    @(36, "ENV", "emi_N2O_f_ind"),
    @(37, "NRG", "pri_coal"),
    @(38, "NRG", "sec_heating_oil"),
    @(39, "NRG", "iip_steel_blafu_slag"),
    @(40, "NRG", "pri_solar_radiation"),
    @(41, "NRG", "sec_hydrogen"),
    @(42, "NRG", "pri_wind_energy_on"),
    @(43, "MAT", "iip_steel_iron_ore"),
    @(44, "NRG", "sec_methane"),
    @(45, "NRG", "pri_deuterium")
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
}
